$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-12-13 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-12-14 Sunday", 2)

$d.Content.Find.Execute("482÷9=53, 5", $true, $false, $false, $false, $false, $true, 1, $false, "244÷6=40, 4", 2)
$d.Content.Find.Execute("322÷4=80, 2", $true, $false, $false, $false, $false, $true, 1, $false, "392÷6=65, 2", 2)
$d.Content.Find.Execute("449÷3=149, 2", $true, $false, $false, $false, $false, $true, 1, $false, "595÷9=66, 1", 2)
$d.Content.Find.Execute("218÷4=54, 2", $true, $false, $false, $false, $false, $true, 1, $false, "280÷7=40, 0", 2)
$d.Content.Find.Execute("873÷5=174, 3", $true, $false, $false, $false, $false, $true, 1, $false, "891÷8=111, 3", 2)

$d.Content.Find.Execute("884÷4=221, 0", $true, $false, $false, $false, $false, $true, 1, $false, "391÷8=48, 7", 2)
$d.Content.Find.Execute("732÷6=122, 0", $true, $false, $false, $false, $false, $true, 1, $false, "349÷5=69, 4", 2)
$d.Content.Find.Execute("827÷7=118, 1", $true, $false, $false, $false, $false, $true, 1, $false, "301÷6=50, 1", 2)
$d.Content.Find.Execute("242÷7=34, 4", $true, $false, $false, $false, $false, $true, 1, $false, "459÷4=114, 3", 2)
$d.Content.Find.Execute("476÷7=68, 0", $true, $false, $false, $false, $false, $true, 1, $false, "369÷7=52, 5", 2)

$d.Content.Find.Execute("666÷7=95, 1", $true, $false, $false, $false, $false, $true, 1, $false, "326÷4=81, 2", 2)
$d.Content.Find.Execute("208÷3=69, 1", $true, $false, $false, $false, $false, $true, 1, $false, "113÷6=18, 5", 2)
$d.Content.Find.Execute("587÷2=293, 1", $true, $false, $false, $false, $false, $true, 1, $false, "503÷4=125, 3", 2)
$d.Content.Find.Execute("241÷8=30, 1", $true, $false, $false, $false, $false, $true, 1, $false, "543÷7=77, 4", 2)
$d.Content.Find.Execute("739÷7=105, 4", $true, $false, $false, $false, $false, $true, 1, $false, "514÷5=102, 4", 2)

$d.Content.Find.Execute("215÷3=71, 2", $true, $false, $false, $false, $false, $true, 1, $false, "779÷2=389, 1", 2)
$d.Content.Find.Execute("226÷5=45, 1", $true, $false, $false, $false, $false, $true, 1, $false, "542÷3=180, 2", 2)
$d.Content.Find.Execute("633÷3=211, 0", $true, $false, $false, $false, $false, $true, 1, $false, "905÷4=226, 1", 2)
$d.Content.Find.Execute("403÷3=134, 1", $true, $false, $false, $false, $false, $true, 1, $false, "729÷5=145, 4", 2)
$d.Content.Find.Execute("388÷4=97, 0", $true, $false, $false, $false, $false, $true, 1, $false, "413÷8=51, 5", 2)

$d.Content.Find.Execute("641÷4=160, 1", $true, $false, $false, $false, $false, $true, 1, $false, "398÷2=199, 0", 2)
$d.Content.Find.Execute("621÷5=124, 1", $true, $false, $false, $false, $false, $true, 1, $false, "594÷2=297, 0", 2)
$d.Content.Find.Execute("930÷5=186, 0", $true, $false, $false, $false, $false, $true, 1, $false, "352÷4=88, 0", 2)
$d.Content.Find.Execute("590÷9=65, 5", $true, $false, $false, $false, $false, $true, 1, $false, "833÷8=104, 1", 2)
$d.Content.Find.Execute("615÷3=205, 0", $true, $false, $false, $false, $false, $true, 1, $false, "808÷4=202, 0", 2)
